$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 102. This pushes the existing rows 102-119
# down to 103-120, carrying along their data/formatting automatically.
$ws.Rows.Item(102).Insert()

# Fill the newly inserted row 102 with the new weekly data point.
# (Columns A, B, C, E, F, G, H, I, J, K, Q, T are constant for this
# market/product subset and match the rows immediately around it.)
$ws.Range("A102").Value = 9
$ws.Range("B102").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C102").Value = "Metropolitana"
$ws.Range("D102").Value = 44946
$ws.Range("E102").Value = 13
$ws.Range("F102").Value = "Fruta"
$ws.Range("G102").Value = 100101
$ws.Range("H102").Value = "Berries"
$ws.Range("I102").Value = 100101004
$ws.Range("J102").Value = "Frambuesa"
$ws.Range("K102").Value = "Sin especificar"
$ws.Range("L102").Value = "Primera"
$ws.Range("M102").Value = 280
$ws.Range("N102").Value = 8000
$ws.Range("O102").Value = 8000
$ws.Range("P102").Value = 8000
$ws.Range("Q102").Value = "$/bandeja 2 kilos"
$ws.Range("R102").Value = "Región de O'Higgins"
$ws.Range("S102").Value = 4000
$ws.Range("T102").Value = 2
